$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "217.87", "0.5267") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.109.77'
$ws.Range("E2").Value = '  -0.11%  '

# Row 3
$ws.Range("D3").Value = '1.655.36'
$ws.Range("E3").Value = '  -0.16%  '

# Row 4
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").Value = '217.87'
$ws.Range("E5").Value = '  +0.59%  '

# Row 6
$ws.Range("D6").Value = '0.5267'
$ws.Range("E6").Value = '  +1.97%  '

# Row 8
$ws.Range("E8").Value = '  -0.93%  '

# Row 9
$ws.Range("D9").Value = '0.06349'
$ws.Range("E9").Value = '  +1.23%  '

# Row 10
$ws.Range("D10").Value = '20.46'
$ws.Range("E10").Value = '  -1.41%  '

# Row 11
$ws.Range("D11").Value = '0.07796'
$ws.Range("E11").Value = '  +1.11%  '

# Row 12
$ws.Range("D12").Value = '4.509'
$ws.Range("E12").Value = '  +1.81%  '

# Row 13
$ws.Range("D13").Value = '1.649.72'
$ws.Range("E13").Value = '  -0.18%  '

# Row 14
$ws.Range("D14").Value = '0.5492'
$ws.Range("E14").Value = '  +1.50%  '

# Row 15
$ws.Range("E15").Value = '  +1.40%  '

# Row 16
$ws.Range("D16").Value = '65.44'
$ws.Range("E16").Value = '  +0.96%  '

# Row 17
$ws.Range("D17").Value = '26.126.86'
$ws.Range("E17").Value = '  -0.12%  '

# Row 18
$ws.Range("E18").Value = '  -0.20%  '

# Row 19
$ws.Range("D19").Value = '4.583'
$ws.Range("E19").Value = '  -0.76%  '

# Row 20
$ws.Range("D20").Value = '190.88'
$ws.Range("E20").Value = '  -0.31%  '

# Row 21
$ws.Range("E21").Value = '  -0.19%  '

# Row 22
$ws.Range("D22").Value = '6.032'
$ws.Range("E22").Value = '  +0.44%  '

# Row 23
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").Value = '141.32'
$ws.Range("E24").Value = '  +1.10%  '

# Row 25
$ws.Range("D25").Value = '0.1236'
$ws.Range("E25").Value = '  +0.99%  '

# Row 26
$ws.Range("D26").Value = '7.252'
$ws.Range("E26").Value = '  +0.90%  '

# Row 27
$ws.Range("D27").Value = '16.09'
$ws.Range("E27").Value = '  +0.11%  '

# Row 28
$ws.Range("D28").Value = '1.427'
$ws.Range("E28").Value = '  +1.12%  '

# Row 29
$ws.Range("D29").Value = '0.05893'
$ws.Range("E29").Value = '  -1.21%  '

# Row 30
$ws.Range("D30").Value = '1.274'
$ws.Range("E30").Value = '  +0.27%  '

# Row 31
$ws.Range("D31").Value = '3.519'
$ws.Range("E31").Value = '  -0.98%  '

# Row 32
$ws.Range("D32").Value = '3.265'
$ws.Range("E32").Value = '  +0.27%  '

# Row 33
$ws.Range("D33").Value = '1.591'
$ws.Range("E33").Value = '  -0.63%  '

# Row 34
$ws.Range("D34").Value = '0.9524'
$ws.Range("E34").Value = '  -1.14%  '

# Row 35
$ws.Range("D35").Value = '2.785'
$ws.Range("E35").Value = '  +0.56%  '

# Row 36
$ws.Range("E36").Value = '  -0.56%  '

# Row 37
$ws.Range("D37").Value = '0.5701'
$ws.Range("E37").Value = '  +0.42%  '

# Row 38
$ws.Range("D38").Value = '0.01617'

# Row 39
$ws.Range("D39").Value = '5.811'
$ws.Range("E39").Value = '  -2.38%  '

# Row 40
$ws.Range("D40").Value = '0.8491'
$ws.Range("E40").Value = '  -0.70%  '

# Row 41
$ws.Range("E41").Value = '  -0.06%  '

# Row 42
$ws.Range("D42").Value = '1.030.25'
$ws.Range("E42").Value = '  +2.38%  '

# Row 43
$ws.Range("D43").Value = '102.62'
$ws.Range("E43").Value = '  +2.03%  '

# Row 44
$ws.Range("D44").Value = '1.799.34'
$ws.Range("E44").Value = '  +0.01%  '

# Row 45
$ws.Range("D45").Value = '57.17'
$ws.Range("E45").Value = '  +0.82%  '

# Row 46
$ws.Range("E46").Value = '  -0.35%  '

# Row 47
$ws.Range("E47").Value = '  +2.85%  '

# Row 48
$ws.Range("E48").Value = '  +2.09%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05151'
$ws.Range("E49").Value = '  -0.34%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.844'
$ws.Range("E50").Value = '  -1.87%  '

# Row 51
$ws.Range("D51").Value = '0.09696'
$ws.Range("E51").Value = '  -0.40%  '

# Restore the default (unstyled) cell style on column D so the saved
# workbook does not retain an explicit Text number format / style index,
# keeping cell styling identical to the original file.
$ws.Range("D2:D51").Style = "Normal"
